$d = $word.ActiveDocument

# --- wdReplaceAll constant ---
$wdReplaceAll = 2

# 1) "...in all aspects CMC corpora..." -> "...in all aspects of CMC corpora..."
#    and "...spoken and and multimodal..." -> "...spoken and multimodal..."
$d.Content.Find.Execute(
    "in all aspects CMC corpora and their evaluation, including spoken and and multimodal interaction.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in all aspects of CMC corpora and their evaluation, including spoken and multimodal interaction.",
    $wdReplaceAll)

# 2) "...up to six (6) pages of content." -> "...up to six (6) pages of content,"
$d.Content.Find.Execute(
    "six (6) pages of content.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "six (6) pages of content,",
    $wdReplaceAll)

# 3) "Papers must be of original, previously-unpublished work." -> "Papers must be original, previously unpublished work."
$d.Content.Find.Execute(
    "Papers must be of original, previously-unpublished work. Papers ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Papers must be original, previously unpublished work. Papers ",
    $wdReplaceAll)

# 4) "omit Sections 1" merge (text already reads correctly; re-assert to normalize run split)
$d.Content.Find.Execute(
    "omit Sections 1",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "omit Sections 1",
    $wdReplaceAll)

# 5) Remove the mid-sentence page break split before "definitive page numbering..."
$d.Content.Find.Execute(
    " Committee will insert the definitive page numbering of articles published in the proceedings",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Committee will insert the definitive page numbering of articles published in the proceedings",
    $wdReplaceAll)

# 6) Merge "should be separated with a semicolon:" + trailing space into one run
$d.Content.Find.Execute(
    "should be separated with a semicolon: (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "should be separated with a semicolon: (",
    $wdReplaceAll)

# 7) Remove the lastRenderedPageBreak split before the "Acknowledgments" heading
$d.Content.Find.Execute(
    "Acknowledgments",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Acknowledgments",
    $wdReplaceAll)

# 8) Merge "papers in conference " + "proceedings, books, journal articles, and book chapters."
$d.Content.Find.Execute(
    "papers in conference proceedings, books, journal articles, and book chapters.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "papers in conference proceedings, books, journal articles, and book chapters.",
    $wdReplaceAll)

# 9) "social Media." -> "social media."
$d.Content.Find.Execute(
    "lexical change in social Media. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "lexical change in social media. ",
    $wdReplaceAll)

# 10) Fix the language tag (de-DE -> en-US) on the PLoS ONE citation pieces
$rng = $d.Content
$rng.Find.Execute("PLoS ONE", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.LanguageID = 1033

$rng2 = $d.Content
$rng2.Find.Execute("PLoS ONE,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Collapse(0)
$rng2.MoveEnd(1, 1)
$rng2.LanguageID = 1033

$rng3 = $d.Content
$rng3.Find.Execute("9(11):e113114.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.LanguageID = 1033

# 11) Italicize "Proceedings of the Eight International Conference on Language Resources and Evaluation (LREC'12)"
$rng4 = $d.Content
$rng4.Find.Execute(
    "Proceedings of the" + [char]0x0020 + "Eight International Conference on Language Resources and Evaluation (LREC" + [char]0x2019 + "12)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng4.Italic = 1

# 12) Add Header/Footer (and linked char) styles by touching the section's header/footer,
#     matching the template's built-in style registration.
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$ftr = $sec.Footers(1)
if ($hdr.Exists) {
    $hdr.Range.Text = $hdr.Range.Text
}
if ($ftr.Exists) {
    $ftr.Range.Text = $ftr.Range.Text
}
